$d = $word.ActiveDocument

# --- Edit 1: after the run "сведения" insert a new run " о программе" ---
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("сведения", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $insStart = $r.End
    $r.Collapse(0)
    $r.InsertAfter(" о программе")

    # Nudge the formatting of just the inserted text so the engine keeps it
    # as its own run instead of silently re-merging it with the neighbour.
    $newRange = $d.Range($insStart, $insStart + 12)
    $newRange.Font.Size = 10
    $newRange.Font.Size = 12
}

# --- Edit 2: tighten wording in the "Для работы с базой данных ..." paragraph ---
$d.Content.Find.Execute("Для работы с базой данных нужно реализовать следующие задачи:", $true, $false, $false, $false, $false, $true, 1, $false, "Для работы с базой данных следующие задачи:", 2)
